# Sprint #2 planning update: add "Sprint 2" section with three new backlog
# items to both the Product Backlog ("Rejestr Produktu") and the Sprint #1
# backlog ("Rejestr Sprint #1") sheets, then move the active tab/selection
# from "Rejestr Sprint #1" to "Rejestr Produktu".

$wb = $excel.ActiveWorkbook

$sprintHeading   = "Sprint 2"
$taskPunctuation = "Jako użytkownik mogę przekształcić tekst tak, aby posiadał znaki interpunkcyjne w odpowiednich miejscach (np. Myślę że tak Albo jednak nie -> Myślę,  że tak. Albo jednak nie.)"
$taskLeetspeak   = "Jako użytkownik mogę zmieniać tekst na leetspeak (leetspeak -> 1337sp34k)"
$taskPolishChars = "Jako użytkownik mogę usunąć polskie znaki z tekstu"

# ---------------------------------------------------------------------
# "Rejestr Produktu" (Product backlog) sheet - append rows 12-16
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Rejestr Produktu")

$ws2.Range("A12").Value = $sprintHeading
$ws2.Range("A12").Font.Bold = $true

$ws2.Range("A13").WrapText = $true
$ws2.Range("A13").VerticalAlignment = -4160

$ws2.Range("A14").Value = $taskPunctuation
$ws2.Range("A14").WrapText = $true
$ws2.Range("A14").VerticalAlignment = -4160
$ws2.Rows.Item(14).RowHeight = 47.25

$ws2.Range("A15").Value = $taskLeetspeak
$ws2.Range("A15").WrapText = $true
$ws2.Range("A15").VerticalAlignment = -4160
$ws2.Rows.Item(15).RowHeight = 31.5

$ws2.Range("A16").Value = $taskPolishChars
$ws2.Range("A16").WrapText = $true
$ws2.Range("A16").VerticalAlignment = -4160

# ---------------------------------------------------------------------
# "Rejestr Sprint #1" sheet - append rows 16, 18-20 (row 17 left blank)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Rejestr Sprint #1")

$ws3.Range("A16").Value = $sprintHeading
$ws3.Range("A16").Font.Bold = $true

$ws3.Range("A18").Value = $taskPunctuation
$ws3.Range("A18").WrapText = $true
$ws3.Range("A18").VerticalAlignment = -4160
$ws3.Rows.Item(18).RowHeight = 63

$ws3.Range("A19").Value = $taskLeetspeak
$ws3.Range("A19").WrapText = $true
$ws3.Range("A19").VerticalAlignment = -4160
$ws3.Rows.Item(19).RowHeight = 31.5

$ws3.Range("A20").Value = $taskPolishChars
$ws3.Range("A20").WrapText = $true
$ws3.Range("A20").VerticalAlignment = -4160

# ---------------------------------------------------------------------
# Move the active tab / selection from "Rejestr Sprint #1" to
# "Rejestr Produktu", and update the selected cells on each sheet.
# ---------------------------------------------------------------------
$ws3.Range("A16:A20").Select()

$ws2.Activate()
$ws2.Range("A21").Select()
